$d = $word.ActiveDocument

function Find-Replace {
    param($findText, $replaceText)
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: find/replace failed for:" $findText
    }
}

function Append-Runs {
    param($paraIndex, $runs)
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Collapse(0)
    foreach ($run in $runs) {
        $text = $run[0]
        $italic = $run[1]
        $runStart = $r.Start
        $r.InsertAfter($text)
        $r.Collapse(0)
        if ($italic) {
            $fmtRange = $d.Range($runStart, $r.Start)
            $fmtRange.Font.Italic = $true
        }
    }
}

function New-BibParagraph {
    param($afterParaIndex, $runs)
    $p = $d.Paragraphs.Item($afterParaIndex)
    $endRange = $p.Range
    $insertPoint = $d.Range($endRange.End, $endRange.End)
    $insertPoint.InsertAfter("`r")
    $insertPoint.Collapse(0)

    $newPara = $d.Paragraphs.Item($afterParaIndex + 1)
    $newPara.Style = "Bibliography"
    $r = $newPara.Range
    $r.MoveEnd(1, -1) | Out-Null

    foreach ($run in $runs) {
        $text = $run[0]
        $italic = $run[1]
        $runStart = $r.Start
        $r.InsertAfter($text)
        $r.Collapse(0)
        if ($italic) {
            $fmtRange = $d.Range($runStart, $r.Start)
            $fmtRange.Font.Italic = $true
        }
    }
}

# ---------------------------------------------------------------------
# 1. Typo fixes in paragraph 6 ("The datasets used ...")
# ---------------------------------------------------------------------
Find-Replace "The datasets used" "The data sets used"
Find-Replace "a variaty of sources" "a variety of sources"

# ---------------------------------------------------------------------
# 2. Typo fixes in paragraph 7 ("... surveys that track changes ...")
# ---------------------------------------------------------------------
Find-Replace "are also caried out" "are also carried out"
Find-Replace "All of the abovementioned datasets" "All of the above-mentioned data sets"
Find-Replace "the volume and viriaty of these data" "the volume and variety of these data"
Find-Replace "the details of each dataset" "the details of each data set"

# ---------------------------------------------------------------------
# 3. Typo fix in paragraph 8 ("Synthesizing data from multiple sources ...")
# ---------------------------------------------------------------------
Find-Replace "information may be lost and uncertinty" "information may be lost and uncertainty"

Write-Host "Done with typo fixes"

# ---------------------------------------------------------------------
# 4. Append the large new passage to the end of paragraph 7 (replacing the
#    trailing "." with an extended discussion). We delete the trailing
#    period first, then append the full continuation (which itself starts
#    with a period).
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$p7r = $p7.Range
$p7r.MoveEnd(1, -1) | Out-Null
$lastChar = $d.Range($p7r.End - 1, $p7r.End)
Write-Host "Last char of para 7:" $lastChar.Text
$lastChar.Delete()

$newTailRuns = New-Object System.Collections.ArrayList
$newTailRuns.Add(@(". Such issues have largely been curtailed in contemporary stock assessments thanks to advances in software that have facilitated the analysis of all available data, in as raw a form as appropriate, in a single integrated analysis", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Maunder and Punt 2013)", $false)) | Out-Null
$newTailRuns.Add(@(". Specifically, the application of statistical modeling tools such as JAGS", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Plummer and others 2003)", $false)) | Out-Null
$newTailRuns.Add(@(", AD Model Builder", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Fournier et al. 2012)", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("and Template Model Builder", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Kristensen et al. 2015)", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("allow the construction of a joint likelihood for an array of observations to, in theory, extract as much information as possible about the biological and fishery processes. However, integrated analyses are not a panacea because model misspecifications and data conflicts are an inevitable consequence of simplifying reality to a small series of equations", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Maunder and Piner 2017)", $false)) | Out-Null
$newTailRuns.Add(@(". A potential solution to this quandary to use a superensemble model, whereby multiple models with different structures are run and their predictions are supplied as covariates to an additional superensemble model", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Anderson et al. 2017)", $false)) | Out-Null
$newTailRuns.Add(@(". Ensemble approaches reduce the risk of picking the", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("“", $false)) | Out-Null
$newTailRuns.Add(@("wrong", $false)) | Out-Null
$newTailRuns.Add(@("”", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("model and also expands the range of hypotheses explored", $false)) | Out-Null
$newTailRuns.Add(@(" ", $false)) | Out-Null
$newTailRuns.Add(@("(Dietterich 2000)", $false)) | Out-Null
$newTailRuns.Add(@(". These advances greatly improve our ability to assess the status and trends of fish populations, however, modern stock assessment biologists are now faced with the overwhelming task of understanding an ever expanding array of data inputs and model outputs.", $false)) | Out-Null

Append-Runs 7 $newTailRuns

Write-Host "Para7 final:" $d.Paragraphs.Item(7).Range.Text

# ---------------------------------------------------------------------
# 5. Delete paragraph 8 entirely ("… Advances in software …")
# ---------------------------------------------------------------------
Write-Host "Para8 (to delete):" $d.Paragraphs.Item(8).Range.Text
$d.Paragraphs.Item(8).Range.Delete()
Write-Host "Para8 after delete:" $d.Paragraphs.Item(8).Range.Text

# ---------------------------------------------------------------------
# 6. Fix the hyphen character in the quote paragraph (now paragraph 8):
#    "data‐intensive" (U+2010 hyphen) -> "data-intensive" (regular hyphen)
# ---------------------------------------------------------------------
Find-Replace "data$([char]0x2010)intensive" "data-intensive"
Write-Host "Para8 after hyphen fix:" $d.Paragraphs.Item(8).Range.Text

# ---------------------------------------------------------------------
# 7. Insert new bibliography entries, in alphabetical order.
#    Current order (after step 6): 9 References heading, 10 Hampton,
#    11 Hilborn, 12 Link, 13 Maunder/Punt 2013, 14 Maunder/Schnute/Ianelli.
# ---------------------------------------------------------------------

# 7a. Anderson et al. 2017 -- inserted right after the References heading (para 9)
$andersonRuns = New-Object System.Collections.ArrayList
$andersonRuns.Add(@("Anderson, Sean C, Andrew B Cooper, Olaf P Jensen, Cóilín Minto, James T Thorson, Jessica C Walsh, Jamie Afflerbach, et al. 2017. “Improving Estimates of Population Status and Trend with Superensemble Models.”", $false)) | Out-Null
$andersonRuns.Add(@(" ", $false)) | Out-Null
$andersonRuns.Add(@("Fish and Fisheries", $true)) | Out-Null
$andersonRuns.Add(@(" ", $false)) | Out-Null
$andersonRuns.Add(@("18 (4). Wiley Online Library: 732–41.", $false)) | Out-Null
New-BibParagraph 9 $andersonRuns

# 7b. Dietterich 2000 -- inserted after Anderson (para 10)
$dietterichRuns = New-Object System.Collections.ArrayList
$dietterichRuns.Add(@("Dietterich, Thomas G. 2000. “Ensemble Methods in Machine Learning.” In", $false)) | Out-Null
$dietterichRuns.Add(@(" ", $false)) | Out-Null
$dietterichRuns.Add(@("International Workshop on Multiple Classifier Systems", $true)) | Out-Null
$dietterichRuns.Add(@(", 1–15. Springer.", $false)) | Out-Null
New-BibParagraph 10 $dietterichRuns

# 7c. Fournier et al. 2012 -- inserted after Dietterich (para 11)
$fournierRuns = New-Object System.Collections.ArrayList
$fournierRuns.Add(@("Fournier, David A, Hans J Skaug, Johnoel Ancheta, James Ianelli, Arni Magnusson, Mark N Maunder, Anders Nielsen, and John Sibert. 2012. “AD Model Builder: Using Automatic Differentiation for Statistical Inference of Highly Parameterized Complex Nonlinear Models.”", $false)) | Out-Null
$fournierRuns.Add(@(" ", $false)) | Out-Null
$fournierRuns.Add(@("Optimization Methods and Software", $true)) | Out-Null
$fournierRuns.Add(@(" ", $false)) | Out-Null
$fournierRuns.Add(@("27 (2). Taylor & Francis: 233–49.", $false)) | Out-Null
New-BibParagraph 11 $fournierRuns

Write-Host "After inserting Anderson/Dietterich/Fournier, count:" $d.Paragraphs.Count

# 7d. Kristensen et al. 2015 -- inserted after Hilborn (para 14), before Link
$kristensenRuns = New-Object System.Collections.ArrayList
$kristensenRuns.Add(@("Kristensen, Kasper, Anders Nielsen, Casper W Berg, Hans Skaug, and Brad Bell. 2015. “TMB: Automatic Differentiation and Laplace Approximation.”", $false)) | Out-Null
$kristensenRuns.Add(@(" ", $false)) | Out-Null
$kristensenRuns.Add(@("arXiv Preprint arXiv:1509.00660", $true)) | Out-Null
$kristensenRuns.Add(@(".", $false)) | Out-Null
New-BibParagraph 14 $kristensenRuns

Write-Host "After inserting Kristensen, count:" $d.Paragraphs.Count
Write-Host "Para15 (should be Kristensen):" $d.Paragraphs.Item(15).Range.Text

# ---------------------------------------------------------------------
# 8. Turn the existing Maunder & Punt 2013 paragraph (now para 17) into the
#    Maunder & Piner 2017 entry by swapping its first and last run text,
#    then insert a brand-new paragraph after it with the original
#    Maunder & Punt 2013 content (so both entries exist, Piner first).
# ---------------------------------------------------------------------
Find-Replace "Maunder, Mark N, and André E Punt. 2013. “A Review of Integrated Analysis in Fisheries Stock Assessment.”" "Maunder, Mark N, and Kevin R Piner. 2017. “Dealing with Data Conflicts in Statistical Inference of Population Assessment Models That Integrate Information from Multiple Diverse Data Sets.”"
Find-Replace "142. Elsevier: 61–74." "192. Elsevier: 16–27."

Write-Host "Para17 (should now be Piner):" $d.Paragraphs.Item(17).Range.Text

$maunderPuntRuns = New-Object System.Collections.ArrayList
$maunderPuntRuns.Add(@("Maunder, Mark N, and André E Punt. 2013. “A Review of Integrated Analysis in Fisheries Stock Assessment.”", $false)) | Out-Null
$maunderPuntRuns.Add(@(" ", $false)) | Out-Null
$maunderPuntRuns.Add(@("Fisheries Research", $true)) | Out-Null
$maunderPuntRuns.Add(@(" ", $false)) | Out-Null
$maunderPuntRuns.Add(@("142. Elsevier: 61–74.", $false)) | Out-Null
New-BibParagraph 17 $maunderPuntRuns

Write-Host "After inserting Maunder/Punt 2013, count:" $d.Paragraphs.Count
Write-Host "Para18 (should be Maunder/Punt 2013):" $d.Paragraphs.Item(18).Range.Text

# ---------------------------------------------------------------------
# 9. Plummer and others 2003 -- appended at the very end, after
#    Maunder/Schnute/Ianelli 2009 (now para 19)
# ---------------------------------------------------------------------
Write-Host "Para19 (should be Maunder/Schnute/Ianelli):" $d.Paragraphs.Item(19).Range.Text

$plummerRuns = New-Object System.Collections.ArrayList
$plummerRuns.Add(@("Plummer, Martyn, and others. 2003. “JAGS: A Program for Analysis of Bayesian Graphical Models Using Gibbs Sampling.” In", $false)) | Out-Null
$plummerRuns.Add(@(" ", $false)) | Out-Null
$plummerRuns.Add(@("Proceedings of the 3rd International Workshop on Distributed Statistical Computing", $true)) | Out-Null
$plummerRuns.Add(@(". Vol. 124. 125.10. Vienna, Austria.", $false)) | Out-Null
New-BibParagraph 19 $plummerRuns

Write-Host "Final paragraph count:" $d.Paragraphs.Count
